$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.073.58'
$ws.Range('E2').Value = '  -3.53%  '
$ws.Range('D3').Value = '1.916.69'
$ws.Range('E3').Value = '  -2.88%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -1.41%  '
$ws.Range('D5').Value = '327.86'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('D7').Value = '0.4688'
$ws.Range('E7').Value = '  -5.72%  '
$ws.Range('D8').Value = '0.4021'
$ws.Range('E8').Value = '  -3.93%  '
$ws.Range('D9').Value = '53.10'
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('D10').Value = '0.08374'
$ws.Range('E10').Value = '  -9.86%  '
$ws.Range('D11').Value = '1.043'
$ws.Range('E11').Value = '  -4.73%  '
$ws.Range('D12').Value = '22.13'
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.888.81'
$ws.Range('E13').Value = '  -5.07%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '7.457'
$ws.Range('E14').Value = '  -5.38%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '6.057'
$ws.Range('E15').Value = '  -6.25%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = '89.57'
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('E18').Value = '  -4.12%  '
$ws.Range('D19').Value = '0.06569'
$ws.Range('E19').Value = '  -2.43%  '
$ws.Range('D20').Value = '17.99'
$ws.Range('E20').Value = '  -5.94%  '
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').Value = '5.713'
$ws.Range('E22').Value = '  -4.15%  '
$ws.Range('D23').Value = '28.053.11'
$ws.Range('E23').Value = '  -3.61%  '
$ws.Range('E24').Value = '  -5.09%  '
$ws.Range('D25').Value = '2.285'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').Value = '2.121.94'
$ws.Range('E26').Value = '  -4.55%  '
$ws.Range('D27').Value = '154.05'
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range('D28').Value = '19.99'
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').Value = '2.131'
$ws.Range('E29').Value = '  -5.90%  '
$ws.Range('D30').Value = '5.691'
$ws.Range('E30').Value = '  -8.59%  '
$ws.Range('D31').Value = '123.04'
$ws.Range('E31').Value = '  -3.31%  '
$ws.Range('D32').Value = '0.9733'
$ws.Range('E32').Value = '  -6.82%  '
$ws.Range('D33').Value = '0.09558'
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('D34').Value = '1.443'
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('D35').Value = '3.634'
$ws.Range('E35').Value = '  -3.07%  '
$ws.Range('D36').Value = '5.529'
$ws.Range('E36').Value = '  -4.80%  '
$ws.Range('D37').Value = '8.801'
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('D38').Value = '0.02302'
$ws.Range('E38').Value = '  -4.73%  '
$ws.Range('D39').Value = '0.06141'
$ws.Range('E39').Value = '  -4.06%  '
$ws.Range('D40').Value = '1.215'
$ws.Range('E40').Value = '  -8.47%  '
$ws.Range('D41').Value = '0.6115'
$ws.Range('E41').Value = '  -5.43%  '
$ws.Range('D42').Value = '11.01'
$ws.Range('E42').Value = '  -4.19%  '
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('D44').Value = '0.1901'
$ws.Range('E44').Value = '  -5.02%  '
$ws.Range('D45').Value = '1.303'
$ws.Range('E45').Value = '  -3.45%  '
$ws.Range('D46').Value = '0.5834'
$ws.Range('E46').Value = '  -5.73%  '
$ws.Range('E47').Value = '  -4.36%  '
$ws.Range('D48').Value = '2.018'
$ws.Range('E48').Value = '  -7.19%  '
$ws.Range('D49').Value = '3.451'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('D50').Value = '0.06842'
$ws.Range('E50').Value = '  -1.72%  '
$ws.Range('D51').Value = '109.39'
$ws.Range('E51').Value = '  -3.32%  '
